$d = $word.ActiveDocument

# In this document, small inline "tag" markers such as
#   <id>p112v_1</id>
# are stored as three separate runs:
#   1) "<id>"    - Courier New, color 7f6000, sz 18 (the opening tag)
#   2) "p112v_1" - plain/black text (the id value)
#   3) "</id>"   - Courier New, color 7f6000, sz 18 (the closing tag)
# This edit collapses each such triplet into a single run carrying the
# opening tag's formatting, with the text being the full concatenation
# "<id>VALUE</id>".

function Merge-IdTag([string]$idValue) {
    $full = "<id>" + $idValue + "</id>"

    $found = $d.Content.Duplicate
    $ok = $found.Find.Execute($full, $false, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
    if (-not $ok) {
        return
    }

    $start = $found.Start
    $end = $found.End
    $openTagLen = 4  # length of "<id>"

    # Range over just the opening "<id>" run - keep its formatting.
    $openRange = $d.Range($start, $start + $openTagLen)
    # Range over the rest ("VALUE</id>") - currently split across two
    # more runs; grab its text then remove it.
    $restRange = $d.Range($start + $openTagLen, $end)
    $restText = $restRange.Text
    $restRange.Delete()

    # Re-append the removed text onto the opening run so everything
    # becomes one run using the opening run's formatting.
    $openRange.InsertAfter($restText)
}

Merge-IdTag "p112v_1"
Merge-IdTag "p113r_1"
